$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Style-preserving donors (kept unchanged elsewhere in the sheet) ---
# C14 = text style (s=14); D15 = numeric style (s=15)

# --- Plain value updates (style unchanged) ---
$ws.Range("M14").Value = -15.384615384615
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -80
$ws.Range("J15").Value = 45
$ws.Range("K15").Value = -15.555555555555
$ws.Range("N15").Value = -42.424242424242
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -23.076923076923
$ws.Range("F16").Value = 41
$ws.Range("G16").Value = 71
$ws.Range("H16").Value = -42.253521126760
$ws.Range("I16").Value = 477
$ws.Range("J16").Value = 667
$ws.Range("K16").Value = -28.485757121439
$ws.Range("L16").Value = 6
$ws.Range("M16").Value = 3.695652173913
$ws.Range("N16").Value = -73.934426229508
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -6.25
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 67
$ws.Range("H17").Value = -20.895522388059
$ws.Range("I17").Value = 754
$ws.Range("J17").Value = 711
$ws.Range("K17").Value = 6.047819971870
$ws.Range("L17").Value = 19.872813990461
$ws.Range("M17").Value = 65.714285714285
$ws.Range("N17").Value = -11.084905660377
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("I18").Value = 289
$ws.Range("J18").Value = 278
$ws.Range("K18").Value = 3.956834532374
$ws.Range("L18").Value = 31.963470319634
$ws.Range("M18").Value = -16.231884057971
$ws.Range("N18").Value = -82.569360675512
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -30.434782608695
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = -4
$ws.Range("I19").Value = 816
$ws.Range("J19").Value = 923
$ws.Range("K19").Value = -11.592632719393
$ws.Range("L19").Value = 19.124087591240
$ws.Range("M19").Value = 60.629921259842
$ws.Range("N19").Value = 15.909090909090
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = 7.142857142857
$ws.Range("F20").Value = 47
$ws.Range("G20").Value = 42
$ws.Range("H20").Value = 11.904761904761
$ws.Range("I20").Value = 581
$ws.Range("J20").Value = 472
$ws.Range("K20").Value = 23.093220338983
$ws.Range("L20").Value = 50.909090909090
$ws.Range("M20").Value = 165.296803652968
$ws.Range("N20").Value = -67.487409065472
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 72
$ws.Range("E21").Value = -15.277777777777
$ws.Range("F21").Value = 244
$ws.Range("G21").Value = 280
$ws.Range("H21").Value = -12.857142857142
$ws.Range("I21").Value = 2966
$ws.Range("J21").Value = 3102
$ws.Range("K21").Value = -4.384268214055
$ws.Range("L21").Value = 22.612649855312
$ws.Range("M21").Value = 45.892769306443
$ws.Range("N21").Value = -57.348288754673
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = -35
$ws.Range("L22").Value = -23.529411764705
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = -70
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 27
$ws.Range("H23").Value = -29.629629629629
$ws.Range("I23").Value = 276
$ws.Range("J23").Value = 299
$ws.Range("K23").Value = -7.692307692307
$ws.Range("L23").Value = 16.455696202531
$ws.Range("M23").Value = 38
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 65
$ws.Range("E24").Value = -52.307692307692
$ws.Range("F24").Value = 137
$ws.Range("G24").Value = 186
$ws.Range("H24").Value = -26.344086021505
$ws.Range("I24").Value = 1725
$ws.Range("J24").Value = 1864
$ws.Range("K24").Value = -7.457081545064
$ws.Range("L24").Value = 37.013502779984
$ws.Range("M24").Value = 29.407351837959
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 95
$ws.Range("H25").Value = -21.052631578947
$ws.Range("I25").Value = 989
$ws.Range("J25").Value = 1008
$ws.Range("K25").Value = -1.884920634920
$ws.Range("L25").Value = 12.004530011325
$ws.Range("M25").Value = -28.437047756874
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 66
$ws.Range("K26").Value = -15.151515151515
$ws.Range("L26").Value = -11.111111111111
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 98
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = 50.769230769230
$ws.Range("L27").Value = 8.888888888888
$ws.Range("I28").Value = 40
$ws.Range("K28").Value = -2.439024390243
$ws.Range("L28").Value = -24.528301886792
$ws.Range("M28").Value = -13.043478260869
$ws.Range("N28").Value = -71.428571428571
$ws.Range("I29").Value = 34
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -22.727272727272
$ws.Range("M29").Value = -12.820512820512
$ws.Range("N29").Value = -73.015873015873

# --- Type/style-changing updates ---
$ws.Range("F14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("D15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("D15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("C29").Value = 2
$ws.Range("D15").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
